$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2276422764227642
$ws.Range("C2").Value = 0.5149051490514905
$ws.Range("J2").Value = 0.002710027100271003
$ws.Range("P2").Value = 0.1544715447154472
$ws.Range("S2").Value = 0.1002710027100271
$ws.Range("C3").Value = 0.0310880829015544
$ws.Range("J3").Value = 0.02072538860103627
$ws.Range("P3").Value = 0.7461139896373057
$ws.Range("S3").Value = 0.2020725388601036
$ws.Range("P4").Value = 0.7692307692307693
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("B6").Value = 0.04265402843601896
$ws.Range("D6").Value = 0.009478672985781991
$ws.Range("F6").Value = 0.06635071090047394
$ws.Range("J6").Value = 0.2796208530805687
$ws.Range("O6").Value = 0.02369668246445497
$ws.Range("Q6").Value = 0.1421800947867299
$ws.Range("R6").Value = 0.05687203791469194
$ws.Range("S6").Value = 0.3791469194312796
$ws.Range("B7").Value = 0.1785714285714286
$ws.Range("D7").Value = 0.02551020408163265
$ws.Range("F7").Value = 0.02551020408163265
$ws.Range("J7").Value = 0.1020408163265306
$ws.Range("O7").Value = 0.02040816326530612
$ws.Range("Q7").Value = 0.1938775510204082
$ws.Range("R7").Value = 0.0663265306122449
$ws.Range("S7").Value = 0.3877551020408163
$ws.Range("B8").Value = 0.101933216168717
$ws.Range("D8").Value = 0.01581722319859402
$ws.Range("E8").Value = 0.00351493848857645
$ws.Range("F8").Value = 0.0421792618629174
$ws.Range("J8").Value = 0.10896309314587
$ws.Range("O8").Value = 0.01757469244288225
$ws.Range("Q8").Value = 0.1652021089630931
$ws.Range("R8").Value = 0.07205623901581722
$ws.Range("S8").Value = 0.4727592267135325
$ws.Range("B9").Value = 0.084070796460177
$ws.Range("D9").Value = 0.01769911504424779
$ws.Range("F9").Value = 0.06637168141592921
$ws.Range("J9").Value = 0.1371681415929203
$ws.Range("Q9").Value = 0.1548672566371681
$ws.Range("R9").Value = 0.1106194690265487
$ws.Range("S9").Value = 0.4292035398230089
$ws.Range("B10").Value = 0.1259968102073365
$ws.Range("D10").Value = 0.01594896331738437
$ws.Range("F10").Value = 0.07177033492822966
$ws.Range("J10").Value = 0.1108452950558214
$ws.Range("O10").Value = 0.01036682615629984
$ws.Range("Q10").Value = 0.1961722488038277
$ws.Range("R10").Value = 0.06778309409888357
$ws.Range("S10").Value = 0.4011164274322169
$ws.Range("G11").Value = 0.167192429022082
$ws.Range("J11").Value = 0.1009463722397476
$ws.Range("K11").Value = 0.2397476340694006
$ws.Range("L11").Value = 0.4700315457413249
$ws.Range("S11").Value = 0.0220820189274448
$ws.Range("G12").Value = 0.8013245033112583
$ws.Range("J12").Value = 0.1721854304635762
$ws.Range("K12").Value = 0.006622516556291391
$ws.Range("L12").Value = 0.01324503311258278
$ws.Range("S12").Value = 0.006622516556291391
$ws.Range("G13").Value = 0.5849056603773585
$ws.Range("J13").Value = 0.2830188679245283
$ws.Range("S13").Value = 0.1320754716981132
$ws.Range("F15").Value = 0.009852216748768473
$ws.Range("H15").Value = 0.1625615763546798
$ws.Range("I15").Value = 0.103448275862069
$ws.Range("J15").Value = 0.3596059113300493
$ws.Range("K15").Value = 0.0541871921182266
$ws.Range("M15").Value = 0.009852216748768473
$ws.Range("O15").Value = 0.02955665024630542
$ws.Range("S15").Value = 0.270935960591133
$ws.Range("F16").Value = 0.02212389380530973
$ws.Range("H16").Value = 0.1858407079646018
$ws.Range("I16").Value = 0.1150442477876106
$ws.Range("J16").Value = 0.3584070796460177
$ws.Range("K16").Value = 0.1017699115044248
$ws.Range("M16").Value = 0.03539823008849557
$ws.Range("O16").Value = 0.03097345132743363
$ws.Range("S16").Value = 0.1504424778761062
$ws.Range("F17").Value = 0.01136363636363636
$ws.Range("H17").Value = 0.2340909090909091
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.3727272727272727
$ws.Range("K17").Value = 0.07272727272727272
$ws.Range("M17").Value = 0.01363636363636364
$ws.Range("O17").Value = 0.05
$ws.Range("S17").Value = 0.1454545454545454
$ws.Range("F18").Value = 0.02272727272727273
$ws.Range("H18").Value = 0.2443181818181818
$ws.Range("I18").Value = 0.07386363636363637
$ws.Range("J18").Value = 0.4034090909090909
$ws.Range("K18").Value = 0.06818181818181818
$ws.Range("M18").Value = 0.01704545454545454
$ws.Range("O18").Value = 0.05113636363636364
$ws.Range("S18").Value = 0.1193181818181818
$ws.Range("F19").Value = 0.01500682128240109
$ws.Range("H19").Value = 0.2380627557980901
$ws.Range("I19").Value = 0.08663028649386084
$ws.Range("J19").Value = 0.3369713506139154
$ws.Range("K19").Value = 0.1064120054570259
$ws.Range("M19").Value = 0.02387448840381992
$ws.Range("O19").Value = 0.06616643929058663
$ws.Range("S19").Value = 0.1268758526603001
